# Script to apply the scraped-data update for
# 2023/algeria_ligue-1_2023-2024.xlsx
#
# 1) Several pairs of existing rows had their match data (columns F:V,
#    i.e. home team .. url_partida) swapped between the two rows of the
#    pair (the Indice/pais/torneio/temporada/data_partida columns A:E
#    stay untouched on each row).
# 2) Two brand-new match rows (87 and 88) are appended at the bottom of
#    the sheet.
#
# NOTE: this runtime's PowerShell-style function calls only bind
# parameters positionally - named arguments (e.g. "-Row 5") silently
# fail to bind. So every helper function below is called positionally.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-MatchRow(
    $Row,
    $Home, $HomeGoals, $Away, $AwayGoals,
    $HomeOpenOdds, $HomeOpenDt, $HomeCloseOdds, $HomeCloseDt,
    $DrawOpenOdds, $DrawOpenDt, $DrawCloseOdds, $DrawCloseDt,
    $AwayOpenOdds, $AwayOpenDt, $AwayCloseOdds, $AwayCloseDt,
    $Url
) {
    $ws.Range("F$Row").Value = $Home
    $ws.Range("G$Row").Value = $HomeGoals
    $ws.Range("H$Row").Value = $Away
    $ws.Range("I$Row").Value = $AwayGoals
    $ws.Range("J$Row").Value = $HomeOpenOdds
    $ws.Range("K$Row").Value = $HomeOpenDt
    $ws.Range("L$Row").Value = $HomeCloseOdds
    $ws.Range("M$Row").Value = $HomeCloseDt
    $ws.Range("N$Row").Value = $DrawOpenOdds
    $ws.Range("O$Row").Value = $DrawOpenDt
    $ws.Range("P$Row").Value = $DrawCloseOdds
    $ws.Range("Q$Row").Value = $DrawCloseDt
    $ws.Range("R$Row").Value = $AwayOpenOdds
    $ws.Range("S$Row").Value = $AwayOpenDt
    $ws.Range("T$Row").Value = $AwayCloseOdds
    $ws.Range("U$Row").Value = $AwayCloseDt
    $ws.Range("V$Row").Value = $Url
}

# --- Pair: rows 21 / 22 -----------------------------------------------
Set-MatchRow 21 "Constantine" 3 "Saoura" 0 `
    1.76 "05/10/2023 07:24" 1.75 "06/10/2023 16:25" `
    3.23 "05/10/2023 07:24" 3.26 "06/10/2023 16:28" `
    4.69 "05/10/2023 07:24" 5.64 "06/10/2023 15:05" `
    "https://www.betexplorer.com/football/algeria/ligue-1/constantine-saoura/2RL4UvSs/"

Set-MatchRow 22 "MC Alger" 5 "ES Setif" 3 `
    1.53 "05/10/2023 07:24" 1.52 "06/10/2023 16:29" `
    3.6 "05/10/2023 07:24" 3.71 "06/10/2023 16:34" `
    6.62 "05/10/2023 07:24" 7.77 "06/10/2023 16:34" `
    "https://www.betexplorer.com/football/algeria/ligue-1/mc-alger-es-setif/OtADSIdf/"

# --- Pair: rows 23 / 24 -----------------------------------------------
Set-MatchRow 23 "Oran" 0 "Magra" 0 `
    1.98 "05/10/2023 07:24" 1.65 "06/10/2023 17:55" `
    2.89 "05/10/2023 07:24" 3.28 "06/10/2023 19:03" `
    3.98 "05/10/2023 07:24" 5.91 "06/10/2023 17:55" `
    "https://www.betexplorer.com/football/algeria/ligue-1/oran-magra/WrVYZ04K/"

Set-MatchRow 24 "Biskra" 0 "Paradou" 5 `
    2.18 "05/10/2023 07:24" 2.03 "06/10/2023 19:34" `
    2.84 "05/10/2023 07:24" 3 "06/10/2023 18:05" `
    3.56 "05/10/2023 07:24" 4.33 "06/10/2023 19:34" `
    "https://www.betexplorer.com/football/algeria/ligue-1/biskra-paradou/hhWUzskE/"

# --- Pair: rows 58 / 59 -----------------------------------------------
Set-MatchRow 58 "Saoura" 0 "ES Setif" 0 `
    1.88 "01/12/2023 05:12" 1.75 "02/12/2023 16:40" `
    3.01 "01/12/2023 05:12" 3.32 "02/12/2023 16:40" `
    4.3 "01/12/2023 05:12" 5.44 "02/12/2023 16:40" `
    "https://www.betexplorer.com/football/algeria/ligue-1/saoura-es-setif/Kb2wh3b1/"

Set-MatchRow 59 "ASO Chlef" 2 "US Souf" 0 `
    1.35 "01/12/2023 05:12" 1.32 "02/12/2023 16:27" `
    4.37 "01/12/2023 05:12" 4.85 "02/12/2023 16:27" `
    8.279999999999999 "01/12/2023 05:12" 11.33 "02/12/2023 16:27" `
    "https://www.betexplorer.com/football/algeria/ligue-1/aso-chlef-us-souf/xzDYhqqe/"

# --- Pair: rows 64 / 65 -----------------------------------------------
Set-MatchRow 64 "Magra" 1 "ASO Chlef" 1 `
    2.38 "08/12/2023 15:43" 2.42 "09/12/2023 14:55" `
    2.82 "08/12/2023 15:43" 2.9 "09/12/2023 14:55" `
    3.47 "08/12/2023 15:43" 3.35 "09/12/2023 14:55" `
    "https://www.betexplorer.com/football/algeria/ligue-1/magra-aso-chlef/CjA5ex5g/"

Set-MatchRow 65 "US Souf" 0 "Saoura" 1 `
    3.28 "08/12/2023 15:43" 3.07 "09/12/2023 13:02" `
    2.91 "08/12/2023 15:43" 2.92 "09/12/2023 13:02" `
    2.41 "08/12/2023 15:43" 2.57 "09/12/2023 14:31" `
    "https://www.betexplorer.com/football/algeria/ligue-1/us-souf-saoura/YP0QjEJO/"

# --- Pair: rows 67 / 68 -----------------------------------------------
Set-MatchRow 67 "Ben Aknoun" 1 "CR Belouizdad" 1 `
    6.41 "14/12/2023 00:12" 5.81 "15/12/2023 15:12" `
    3.37 "14/12/2023 00:12" 3.59 "15/12/2023 15:12" `
    1.56 "14/12/2023 00:12" 1.65 "15/12/2023 15:12" `
    "https://www.betexplorer.com/football/algeria/ligue-1/es-ben-aknoun-cr-belouizdad/CjWzBH4t/"

Set-MatchRow 68 "El Bayadh" 1 "USM Alger" 1 `
    2.12 "14/12/2023 01:12" 2.16 "15/12/2023 15:13" `
    2.86 "14/12/2023 01:12" 2.85 "15/12/2023 15:10" `
    3.99 "14/12/2023 01:12" 4.3 "15/12/2023 15:13" `
    "https://www.betexplorer.com/football/algeria/ligue-1/el-bayadh-usm-alger/Ywta5DYH/"

# --- Pair: rows 80 / 81 -----------------------------------------------
Set-MatchRow 80 "Biskra" 1 "Ben Aknoun" 1 `
    1.57 "28/12/2023 07:12" 1.41 "29/12/2023 17:57" `
    3.47 "28/12/2023 07:12" 4.05 "29/12/2023 17:57" `
    5.97 "28/12/2023 07:12" 9.98 "29/12/2023 17:57" `
    "https://www.betexplorer.com/football/algeria/ligue-1/biskra-es-ben-aknoun/z72UKhY4/"

Set-MatchRow 81 "USM Alger" 0 "MC Alger" 0 `
    2.92 "28/12/2023 07:12" 3.62 "29/12/2023 17:59" `
    2.96 "28/12/2023 07:12" 2.98 "29/12/2023 17:56" `
    2.42 "28/12/2023 07:12" 2.24 "29/12/2023 17:59" `
    "https://www.betexplorer.com/football/algeria/ligue-1/usm-alger-mc-alger/jDDqjB3f/"

# --- New rows appended at the bottom: 87 and 88 ------------------------
function Add-MatchRow(
    $Row, $Indice, $DataPartida,
    $Home, $HomeGoals, $Away, $AwayGoals,
    $HomeOpenOdds, $HomeOpenDt, $HomeCloseOdds, $HomeCloseDt,
    $DrawOpenOdds, $DrawOpenDt, $DrawCloseOdds, $DrawCloseDt,
    $AwayOpenOdds, $AwayOpenDt, $AwayCloseOdds, $AwayCloseDt,
    $Url
) {
    $ws.Range("A$Row").Value = $Indice

    $ws.Range("B$Row").Value = "algeria"
    $ws.Range("C$Row").Value = "ligue-1"
    $ws.Range("D$Row").Value = "2023-2024"

    # Column E ("data_partida") keeps the same date/time number format used
    # throughout the sheet (cellXfs index 2 in the original file).
    $ws.Range("E$Row").Value = $DataPartida
    $ws.Range("E$Row").NumberFormat = "YYYY-MM-DD HH:MM:SS"

    Set-MatchRow $Row $Home $HomeGoals $Away $AwayGoals `
        $HomeOpenOdds $HomeOpenDt $HomeCloseOdds $HomeCloseDt `
        $DrawOpenOdds $DrawOpenDt $DrawCloseOdds $DrawCloseDt `
        $AwayOpenOdds $AwayOpenDt $AwayCloseOdds $AwayCloseDt `
        $Url
}

Add-MatchRow 87 86 45297.66666666666 "Saoura" 0 "Khenchela" 1 `
    1.78 "04/01/2024 09:12" 1.52 "06/01/2024 15:57" `
    3.1 "04/01/2024 09:12" 3.62 "06/01/2024 15:57" `
    4.87 "04/01/2024 09:12" 8.130000000000001 "06/01/2024 15:57" `
    "https://www.betexplorer.com/football/algeria/ligue-1/saoura-khenchela/Ec1YJCmB/"

Add-MatchRow 88 87 45297.75 "MC Alger" 1 "Kabylie" 1 `
    1.49 "04/01/2024 09:12" 1.48 "06/01/2024 17:58" `
    3.58 "04/01/2024 09:12" 3.79 "06/01/2024 17:58" `
    7.24 "04/01/2024 09:12" 8.69 "06/01/2024 17:58" `
    "https://www.betexplorer.com/football/algeria/ligue-1/mc-alger-kabylie/4r85CUuo/"

# Column A ("Indice") uses the same bold / bordered / centered look as
# every other row in the sheet (cellXfs index 1 in the original file).
# Build that style once on row 87 and then just copy the formatting
# (not the value) down to row 88, so we do not create throw-away
# intermediate cell styles in the process.
$ws.Range("A87").Font.Bold = $true
$ws.Range("A87").HorizontalAlignment = -4108
$ws.Range("A87").VerticalAlignment = -4160
$ws.Range("A87").Borders.LineStyle = 1

$ws.Range("A87").Copy()
$ws.Range("A88").PasteSpecial(-4122)
